$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Remove the leading "Carátula" and "Índice" paragraphs (each ends in
#    a manual page break) at the very start of the document.
# ---------------------------------------------------------------------
$start = $d.Paragraphs.Item(1).Range.Start
$end = $d.Paragraphs.Item(2).Range.End
$r = $d.Range($start, $end)
$r.Delete()

# ---------------------------------------------------------------------
# 2) Remove the three paragraphs right after the "1.5. Estrategia
#    metodológica" heading:
#       "Se identifican las principales etapas..."
#       "1.- Proceso para la obtención de datos..."
#       "2.- Herramientas, técnicas y/o modelos..."
# ---------------------------------------------------------------------
$pHeading = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*Estrategia metodol*") {
        $pHeading = $i
    }
}

$delStart = $d.Paragraphs.Item($pHeading + 1).Range.Start
$delEnd = $d.Paragraphs.Item($pHeading + 3).Range.End
$r2 = $d.Range($delStart, $delEnd)
$r2.Delete()

# ---------------------------------------------------------------------
# 3) Extend the final paragraph's sentence about the Atlite model with
#    additional detail.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Una vez que se tengan estos datos se introducirán al modelo para simulación que se tiene disponible, desarrollado por Python for Power System Analysis (PyPSA) y llamado Atlite.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Una vez que se tengan estos datos se introducirán al modelo para simulación que se tiene disponible, desarrollado por Python for Power System Analysis (PyPSA) y llamado Atlite. El modelo permitirá simular la producción de energía solar por hora durante el periodo establecido y en la zonas indicadas (para esto último se utilizará la librería GeoPandas de Python); igualmente permitirá seleccionar hacer simulaciones con paneles solares de Silicon o de Cadmio Telurio.",
    2) | Out-Null

# ---------------------------------------------------------------------
# 4) Append the new paragraphs at the end of the document:
#      - closing remark + page break
#      - "2. Desarrollo" heading (bold)
#      - methodology-application paragraph (blue color)
#      - "llllll" paragraph
# ---------------------------------------------------------------------
$last = $d.Paragraphs.Item($d.Paragraphs.Count)
$rng = $last.Range
$rng.Collapse(0)
$rng.InsertAfter("`r" + "Tras realizar las simulaciones, se espera tener suficiente información para analizar la viabilidad energética de los paneles solares flotantes en Alemania y España, de lo cuál se podría inferir la viabilidad de estos para México.")

$last = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastRange = $last.Range
$lastRange.Collapse(0)
$lastRange.InsertBefore("`u{0007}")
# replace trailing marker with actual page break via Find/Replace on the character
$last.Range.Characters.Last.Text = ""

$rng2 = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$rng2.Collapse(0)
$rng2.InsertAfter("`r2. Desarrollo")

$rng3 = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$rng3.Collapse(0)
$rng3.InsertAfter("`rAplicación de la estrategia metodológica describiendo básicamente la investigación de campo realizada, el análisis de los datos, las pruebas o validación de las hipótesis, la discusión de los resultados y las respuestas a las preguntas de investigación.")

$rng4 = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$rng4.Collapse(0)
$rng4.InsertAfter("`rllllll")
